$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-12: only the cells whose contents actually change ---

# Row 2
$ws.Range("B2").Value = "NSE:ADFFOODS"
$ws.Range("C2").Value = "NSE:ABAN"
$ws.Range("E2").Value = "NSE:BAJFINANCE"
$ws.Range("F2").Value = ""

# Row 3
$ws.Range("B3").Value = "NSE:ASIANHOTNR"
$ws.Range("C3").Value = "NSE:AKSHOPTFBR"
$ws.Range("E3").Value = "NSE:BHARTIARTL"

# Row 4
$ws.Range("B4").Value = "NSE:CLSEL"
$ws.Range("C4").Value = "NSE:AMRUTANJAN"
$ws.Range("E4").Value = "NSE:JKCEMENT"

# Row 5
$ws.Range("B5").Value = "NSE:EIHOTEL"
$ws.Range("C5").Value = "NSE:ASAL"
$ws.Range("E5").Value = "NSE:OBEROIRLTY"

# Row 6
$ws.Range("B6").Value = "NSE:FAZE3Q"
$ws.Range("C6").Value = "NSE:ASTRAZEN"
$ws.Range("E6").Value = "NSE:PAGEIND"

# Row 7
$ws.Range("B7").Value = "NSE:GILLANDERS"
$ws.Range("C7").Value = "NSE:CORDSCABLE"

# Row 8
$ws.Range("B8").Value = "NSE:GROBTEA"
$ws.Range("C8").Value = "NSE:ESCORTS"

# Row 9
$ws.Range("B9").Value = "NSE:HINDOILEXP"
$ws.Range("C9").Value = "NSE:GOYALALUM"

# Row 10
$ws.Range("B10").Value = "NSE:JKPAPER"
$ws.Range("C10").Value = "NSE:HCL-INSYS"

# Row 11
$ws.Range("B11").Value = "NSE:LALPATHLAB"
$ws.Range("C11").Value = "NSE:IGPL"

# Row 12
$ws.Range("B12").Value = "NSE:PRIVISCL"
$ws.Range("C12").Value = "NSE:KIRLOSBROS"

# --- Append new rows 13-19 ---
# Column A carries the same bold/centered/bordered style as the existing index column (copy format from A2).
$srcA = $ws.Range("A2")

$newRows = @(
    @{ Row = 13; Idx = 11; B = "NSE:SAFARI";     C = "NSE:KRITI" },
    @{ Row = 14; Idx = 12; B = "NSE:SAKHTISUG";  C = "NSE:LICNFNHGP" },
    @{ Row = 15; Idx = 13; B = $null;            C = "NSE:MANAKSIA" },
    @{ Row = 16; Idx = 14; B = $null;            C = "NSE:MANGALAM" },
    @{ Row = 17; Idx = 15; B = $null;            C = "NSE:PALREDTEC" },
    @{ Row = 18; Idx = 16; B = $null;            C = "NSE:PIXTRANS" },
    @{ Row = 19; Idx = 17; B = $null;            C = "NSE:RENUKA" }
)

foreach ($entry in $newRows) {
    $r = $entry.Row

    $srcA.Copy()
    $ws.Range("A$r").PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = $entry.Idx

    if ($entry.B) {
        $ws.Cells.Item($r, 2).Value = $entry.B
    }
    if ($entry.C) {
        $ws.Cells.Item($r, 3).Value = $entry.C
    }
}
